$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '88.581.31'
$ws.Range("E2").Value = '  +10.37%  '

# Row 3
$ws.Range("D3").Value = '3.360.05'
$ws.Range("E3").Value = '  +5.42%  '

# Row 4
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").Value = "'222.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.02%  '

# Row 6
$ws.Range("D6").Value = "'655.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.48%  '

# Row 7
$ws.Range("D7").Value = "'0.345"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +25.40%  '

# Row 8
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.08%  '

# Row 9
$ws.Range("D9").Value = "'0.619"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.26%  '

# Row 10
$ws.Range("D10").Value = '3.353.17'
$ws.Range("E10").Value = '  +5.32%  '

# Row 11
$ws.Range("D11").Value = "'0.615"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.85%  '

# Row 12
$ws.Range("D12").Value = "'0.0000276"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +6.53%  '

# Row 13
$ws.Range("E13").Value = '  +2.30%  '

# Row 14
$ws.Range("D14").Value = "'35.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +11.39%  '

# Row 15
$ws.Range("D15").Value = '3.979.19'
$ws.Range("E15").Value = '  +5.46%  '

# Row 16
$ws.Range("E16").Value = '  +4.05%  '

# Row 17
$ws.Range("D17").Value = '87.947.44'
$ws.Range("E17").Value = '  +9.55%  '

# Row 18
$ws.Range("D18").Value = '3.352.64'
$ws.Range("E18").Value = '  +5.18%  '

# Row 19
$ws.Range("D19").Value = "'14.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.96%  '

# Row 20
$ws.Range("D20").Value = "'3.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.01%  '

# Row 21
$ws.Range("D21").Value = "'470.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.74%  '

# Row 22
$ws.Range("B22").Value = 'Polkadot'
$ws.Range("C22").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D22").Value = "'5.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +9.57%  '

# Row 23
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = "'9.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.20%  '

# Row 24
$ws.Range("D24").Value = "'13.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +24.45%  '

# Row 25
$ws.Range("E25").Value = '  +7.26%  '

# Row 26
$ws.Range("D26").Value = "'5.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +17.99%  '

# Row 27
$ws.Range("D27").Value = '3.507.90'
$ws.Range("E27").Value = '  +4.78%  '

# Row 28
$ws.Range("D28").Value = "'79.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.52%  '

# Row 29
$ws.Range("E29").Value = '  +63.31%  '

# Row 30
$ws.Range("E30").Value = '  +5.64%  '

# Row 31
$ws.Range("D31").Value = "'0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.13%  '

# Row 32
$ws.Range("E32").Value = '  +5.09%  '

# Row 33
$ws.Range("D33").Value = "'605.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.06%  '

# Row 34
$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").Value = "'1.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.54%  '

# Row 35
$ws.Range("B35").Value = 'Binance-PegBSC-USD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D35").Value = "'0.990"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.11%  '

# Row 36
$ws.Range("D36").Value = "'2.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.57%  '

# Row 37
$ws.Range("E37").Value = '  +0.48%  '

# Row 38
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D38").Value = "'7.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +24.90%  '

# Row 39
$ws.Range("B39").Value = 'EthereumClassic'
$ws.Range("C39").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D39").Value = "'24.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.14%  '

# Row 40
$ws.Range("D40").Value = "'0.425"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.41%  '

# Row 41
$ws.Range("D41").Value = "'2.16"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +19.28%  '

# Row 42
$ws.Range("D42").Value = "'0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.03%  '

# Row 43
$ws.Range("D43").Value = "'21.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.39%  '

# Row 44
$ws.Range("D44").Value = "'3.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +13.69%  '

# Row 45
$ws.Range("D45").Value = "'193.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.89%  '

# Row 46
$ws.Range("E46").Value = '  +0.00%  '

# Row 47
$ws.Range("D47").Value = "'158.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.89%  '

# Row 48
$ws.Range("D48").Value = "'47.81"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +11.91%  '

# Row 49
$ws.Range("D49").Value = "'1.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.75%  '

# Row 50
$ws.Range("D50").Value = "'0.800"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.13%  '

# Row 51
$ws.Range("D51").Value = "'26.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.06%  '
